$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new test-log entry as row 5 (the sheet's used range grows from
# A1:E4 to A1:E5). Columns: Test ID, Test By, Expected Result, Result, Comments.
$ws.Range("A5").Value = "LogoLink"
$ws.Range("B5").Value = "Dan Brown"
$ws.Range("C5").Value = "Opens a webpage"
$ws.Range("D5").Value = "Error 404 no webpage is found"
$ws.Range("E5").Value = "The logo is the  being clicked is the cool little batman"

# Match the wrap-text formatting used on the Expected Result / Result /
# Comments columns for every other row in the table.
$ws.Range("C5:E5").WrapText = $true

# Row grows to fit the wrapped text, same as the other multi-line rows.
$ws.Rows.Item(5).RowHeight = 43.2

$ws.Range("E5").Select() | Out-Null

$wb.Save()
